# Trade #96 closed at 2026-02-17 15:57:18 - unknown UNKNOWN +0.000%
#
# Updates the Summary and Strategy Status roll-up figures for the
# MarketMaking strategy, and appends the newly closed trade (#96) as
# row 97 on both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.92   # Current Capital
$summary.Range("B4").Value = -0.09     # Total P&L $
$summary.Range("B6").Value = 96        # Total Trades
$summary.Range("B7").Value = 36        # Winning Trades
$summary.Range("B9").Value = 37.5      # Win Rate %

# --- Strategy Status sheet (MarketMaking is row 4) --------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.92   # Capital
$status.Range("D4").Value = 96      # Trades
$status.Range("E4").Value = -0.09   # P&L $
$status.Range("F4").Value = -0.08   # P&L %
$status.Range("G4").Value = 37.5    # Win Rate %

# --- Append the newly closed trade as row 97 --------------------------
# Column B holds a literal date-like string ("2026-02-17"), so it is
# entered with a leading quote to keep Excel from auto-converting it to
# a date serial number, matching every other row in the sheet.
$newRowValues = @(
    96,
    "'2026-02-17",
    "15:57:11",
    "MarketMaking",
    "DOWN",
    0.82,
    0.83,
    "CLOSED",
    1.2195,
    0.01,
    99.92,
    0,
    0,
    0.6,
    "Normal spread capture: 19600 bps",
    "early_exit",
    0.13
)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($col = 1; $col -le $newRowValues.Length; $col++) {
        $ws.Cells.Item(97, $col).Value = $newRowValues[$col - 1]
    }
}
